# Update countries & provincias Spain
# Applies the COVID data refresh captured in the commit diff:
#  - bump the "Datos actualizados" timestamp string
#  - update case counters for several countries (rows 4, 7, 13, 50, 136, 140)
#  - re-sort four low-count territories (Santa Sede / Islas Turcas y Caicos /
#    Seychelles / Montserrat), which also shuffles their active/deaths figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 23:11"

# --- Country counter updates -------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 2179892
$ws.Range("C4").Value = 17664
$ws.Range("D4").Value = 878973
$ws.Range("E4").Value = 1182699
$ws.Range("G4").Value = 362
$ws.Range("H4").Value = 118220

# Row 7: India
$ws.Range("B7").Value = 343026
$ws.Range("C7").Value = 10243
$ws.Range("D7").Value = 180320
$ws.Range("E7").Value = 152791
$ws.Range("G7").Value = 395
$ws.Range("H7").Value = 9915

# Row 13: Alemania
$ws.Range("B13").Value = 188044
$ws.Range("C13").Value = 373
$ws.Range("E13").Value = 6559
$ws.Range("G13").Value = 15
$ws.Range("H13").Value = 8885

# Row 50: Barein
$ws.Range("B50").Value = 19013
$ws.Range("C50").Value = 786
$ws.Range("D50").Value = 13267
$ws.Range("E50").Value = 5700

# Row 136: Cabo Verde
$ws.Range("B136").Value = 760
$ws.Range("C136").Value = 10
$ws.Range("D136").Value = 340
$ws.Range("E136").Value = 413
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 7

# Row 140: Santo Tome y Principe
$ws.Range("B140").Value = 662
$ws.Range("C140").Value = 1
$ws.Range("E140").Value = 473

# --- Reorder the four small territories (208-211) ----------------------
# Old order: Santa Sede, Islas Turcas y Caicos, Seychelles, Montserrat
# New order: Islas Turcas y Caicos, Santa Sede, Montserrat, Seychelles
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
